# amazon.xlsx - "added loadProperties into FileUtils"
#
# The RUN column (A) for the second Sign-in/Sign-out test block (rows 13-19)
# is switched from "Y"/"y" to "n" (disabling those rows), row 9's RUN flag
# is normalised from "Y" to "y", a couple of row heights are nudged, the
# active selection moves to A13, and the stored "display" text for the two
# sign-in hyperlinks is updated to include the password that is already
# present in the cell text itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- RUN column updates ------------------------------------------------
$ws.Range("A9").Value  = "y"
$ws.Range("A13").Value = "n"
$ws.Range("A14").Value = "n"
$ws.Range("A15").Value = "n"
$ws.Range("A16").Value = "n"
$ws.Range("A17").Value = "n"
$ws.Range("A18").Value = "n"
$ws.Range("A19").Value = "n"

# --- row height tweaks ---------------------------------------------------
$ws.Rows.Item(9).RowHeight  = 24.65
$ws.Rows.Item(10).RowHeight = 24.65
$ws.Rows.Item(14).RowHeight = 14.35

# --- hyperlink display text ----------------------------------------------
$ws.Range("C5").Hyperlinks.Item(1).TextToDisplay  = "gpawel17@email.com|1Qazxsw2!"
$ws.Range("C14").Hyperlinks.Item(1).TextToDisplay = "gpawel17@email.com|1Qazxsw2!"

# --- selection -------------------------------------------------------------
$ws.Range("A13").Select()
